$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.372.05'
$ws.Range("E2").Value = '  -7.40%  '
$ws.Range("D3").Value = '1.680.77'
$ws.Range("E3").Value = '  -6.15%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.21'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5074'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -13.65%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.006'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2668'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '22.10'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06320'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07378'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.34%  '
$ws.Range("D12").Value = '1.683.03'
$ws.Range("E12").Value = '  -6.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.525'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5791'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.82%  '
$ws.Range("D15").Value = '1.910.10'
$ws.Range("E15").Value = '  -6.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008653'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.48%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.19'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -13.63%  '
$ws.Range("D18").Value = '26.428.31'
$ws.Range("E18").Value = '  -7.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.002'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -7.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.005'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.87'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '186.95'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -10.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.251'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.005'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.71'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.494'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1170'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.91%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.97'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.345'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05773'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.01%  '
$ws.Range("E31").Value = '  -6.11%  '
$ws.Range("E32").Value = '  -6.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.517'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.664'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.011'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.41%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.5958'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.361'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.65%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.676'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.98%  '
$ws.Range("D39").Value = '1.100.31'
$ws.Range("E39").Value = '  -4.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01604'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.898'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.53%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8599'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.64%  '
$ws.Range("E43").Value = '  -0.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.96'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.64%  '
$ws.Range("D45").Value = '1.837.87'
$ws.Range("E45").Value = '  -5.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000116'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.42'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.99%  '
$ws.Range("E48").Value = '  +0.67%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.979'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.77%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4315'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.47%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05215'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.98%  '
